# Update the cached "datetimeFigureOut" date field text from 7/23/2017 to
# 8/12/2017 everywhere it appears: the slide master and every slide layout.

$p = $ppt.ActivePresentation

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "7/23/2017") {
                $tr.Text = "8/12/2017"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShapes $layout.Shapes
}
